# feat: allow creating notes and vocabulary entries with tags
#
# 1) Fix a pre-existing typo in the ENGLISH vocabulary list.
# 2) Append 10 new vocabulary entries to the ENGLISH sheet.
# 3) Append 1 new tagged note to the NOTES sheet.
# 4) Restore selection/active-sheet state to match the author's last edit
#    (ENGLISH sheet selected at M24; NOTES sheet selection back at A1).

$wb = $excel.ActiveWorkbook

$wsEnglish = $wb.Worksheets.Item("ENGLISH")
$wsNotes = $wb.Worksheets.Item("NOTES")

# --- 1) Typo fix -----------------------------------------------------------
$wsEnglish.Range("A67").Value = "disentangle"

# --- 2) New vocabulary entries (ENGLISH sheet, rows 115-124) ---------------
$wsEnglish.Range("A115").Value = "resentment"
$wsEnglish.Range("C115").Value = "bitterness"
$wsEnglish.Range("E115").Value = 0
$wsEnglish.Range("F115").Value = "2021-11-18 13:38:55.766952"

$wsEnglish.Range("A116").Value = "coalesce"
$wsEnglish.Range("B116").Value = "com together to form one mass or whole"
$wsEnglish.Range("C116").Value = "unite"
$wsEnglish.Range("E116").Value = 0
$wsEnglish.Range("F116").Value = "2021-11-18 13:39:52.321676"

$wsEnglish.Range("A117").Value = "timid"
$wsEnglish.Range("C117").Value = "easily frightened"
$wsEnglish.Range("E117").Value = 0
$wsEnglish.Range("F117").Value = "2021-11-18 13:40:17.759858"

$wsEnglish.Range("A118").Value = "reverberation"
$wsEnglish.Range("B118").Value = "prolongation of a sound; a continuing effect"
$wsEnglish.Range("C118").Value = "resonance"
$wsEnglish.Range("E118").Value = 0
$wsEnglish.Range("F118").Value = "2021-11-18 13:41:26.98325"

$wsEnglish.Range("A119").Value = "imprisonment"
$wsEnglish.Range("C119").Value = "incarceration"
$wsEnglish.Range("E119").Value = 0
$wsEnglish.Range("F119").Value = "2021-11-18 13:42:08.650617"

$wsEnglish.Range("A120").Value = "convention"
$wsEnglish.Range("C120").Value = "agreement;custom"
$wsEnglish.Range("E120").Value = 0
$wsEnglish.Range("F120").Value = "2021-11-18 13:44:21.614304"

$wsEnglish.Range("A121").Value = "crestfallen"
$wsEnglish.Range("C121").Value = "disappointed;downhearted"
$wsEnglish.Range("E121").Value = 0
$wsEnglish.Range("F121").Value = "2021-11-18 13:45:08.281708"

$wsEnglish.Range("A122").Value = "innate"
$wsEnglish.Range("C122").Value = "natural;inborn"
$wsEnglish.Range("E122").Value = 0
$wsEnglish.Range("F122").Value = "2021-11-18 13:46:56.18276"

$wsEnglish.Range("A123").Value = "muddle"
$wsEnglish.Range("C123").Value = "confuse;bewilder"
$wsEnglish.Range("E123").Value = 0
$wsEnglish.Range("F123").Value = "2021-11-18 13:47:26.716003"

$wsEnglish.Range("A124").Value = "resolutely"
$wsEnglish.Range("B124").Value = "in an admirably purposeful, determined, and unwavering manner"
$wsEnglish.Range("E124").Value = 0
$wsEnglish.Range("F124").Value = "2021-11-18 13:48:20.311353"

# --- 3) New tagged note (NOTES sheet, row 27) -------------------------------
$wsNotes.Range("A27").Value = "The main thing is to keep the main thing the main thing"
$wsNotes.Range("B27").Value = "essentialism"

# --- 4) Selection / active sheet state --------------------------------------
# Select NOTES!A1 first, then ENGLISH!M24 last, so ENGLISH ends up the
# active tab (index 0) while NOTES keeps a clean A1 selection for next time.
$wsNotes.Range("A1").Select()
$wsEnglish.Range("M24").Select()
